# save data done + era data updated
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1, matching the style of the other headers (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

$h = @(0, 1, 0, 1, 0, 0, 0, 0, 0, 1, 1, 0, 0)
for ($i = 0; $i -lt $h.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $h[$i]
}
